# Updates cached market-price / profit values on the Kraken_Profits sheets
# (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR) following the scheduled runner's
# refreshed pricing snapshot. For each affected leve row, the
# currentAveragePrice* / LevePrice* / LeveProfit* columns (H-N) are
# rewritten with newly computed figures; cells that no longer have a
# value in the refreshed snapshot are cleared rather than zeroed so the
# cell disappears entirely (matching rows where, e.g., the HQ profit
# column had no computable value).

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")

$ws.Range("H2").Value = 390.22223
$ws.Range("I2").Value = 389
$ws.Range("J2").Value = 400
$ws.Range("K2").Value = 389
$ws.Range("L2").Value = 400
$ws.Range("M2").Value = -276
$ws.Range("N2").Value = -626
$ws.Range("H69").Value = 3333
$ws.Range("I69").Value = 3333
$ws.Range("K69").Value = 9999
$ws.Range("M69").Value = -9125
$ws.Range("H72").Value = 3333
$ws.Range("I72").Value = 3333
$ws.Range("K72").Value = 29997
$ws.Range("M72").Value = -25629
$ws.Range("H86").Value = 2000
$ws.Range("J86").Value = 0
$ws.Range("L86").Value = 0
$ws.Range("N86").ClearContents()
$ws.Range("H88").Value = 3840
$ws.Range("I88").Value = 4050
$ws.Range("J88").Value = 3000
$ws.Range("K88").Value = 4050
$ws.Range("L88").Value = 3000
$ws.Range("M88").Value = -3644
$ws.Range("N88").Value = -3812
$ws.Range("H89").Value = 2000
$ws.Range("J89").Value = 0
$ws.Range("L89").Value = 0
$ws.Range("N89").ClearContents()
$ws.Range("H91").Value = 3840
$ws.Range("I91").Value = 4050
$ws.Range("J91").Value = 3000
$ws.Range("K91").Value = 4050
$ws.Range("L91").Value = 3000
$ws.Range("M91").Value = -2646
$ws.Range("N91").Value = -5808
$ws.Range("H141").Value = 595.1667
$ws.Range("I141").Value = 595.1667
$ws.Range("K141").Value = 1785.5001
$ws.Range("M141").Value = 3394.4999

$ws = $wb.Worksheets.Item("ARM")

$ws.Range("H2").Value = 5266
$ws.Range("I2").Value = 5000
$ws.Range("J2").Value = 5399
$ws.Range("K2").Value = 5000
$ws.Range("L2").Value = 5399
$ws.Range("M2").Value = -4887
$ws.Range("N2").Value = -5625
$ws.Range("H3").Value = 1248.6666
$ws.Range("I3").Value = 623
$ws.Range("K3").Value = 623
$ws.Range("M3").Value = -508
$ws.Range("H7").Value = 0
$ws.Range("I7").Value = 0
$ws.Range("K7").Value = 0
$ws.Range("M7").ClearContents()
$ws.Range("H45").Value = 2775
$ws.Range("I45").Value = 0
$ws.Range("K45").Value = 0
$ws.Range("M45").ClearContents()
$ws.Range("H116").Value = 5266
$ws.Range("I116").Value = 5000
$ws.Range("J116").Value = 5399
$ws.Range("K116").Value = 5000
$ws.Range("L116").Value = 5399
$ws.Range("M116").Value = -2706
$ws.Range("N116").Value = -9987
$ws.Range("H132").Value = 5099.8335
$ws.Range("I132").Value = 4569.8
$ws.Range("K132").Value = 13709.4
$ws.Range("M132").Value = -11179.4

$ws = $wb.Worksheets.Item("BSM")

$ws.Range("H3").Value = 5266
$ws.Range("I3").Value = 5000
$ws.Range("J3").Value = 5399
$ws.Range("K3").Value = 5000
$ws.Range("L3").Value = 5399
$ws.Range("M3").Value = -4886
$ws.Range("N3").Value = -5627
$ws.Range("H7").Value = 10000076
$ws.Range("I7").Value = 10000076
$ws.Range("K7").Value = 10000076
$ws.Range("M7").Value = -9999963
$ws.Range("H137").Value = 55000
$ws.Range("I137").Value = 55000
$ws.Range("K137").Value = 55000
$ws.Range("M137").Value = -49900

$ws = $wb.Worksheets.Item("CRP")

$ws.Range("H5").Value = 1320.6666
$ws.Range("I5").Value = 574.6667
$ws.Range("J5").Value = 2066.6667
$ws.Range("K5").Value = 574.6667
$ws.Range("L5").Value = 2066.6667
$ws.Range("M5").Value = -462.6667
$ws.Range("N5").Value = -2290.6667
$ws.Range("H7").Value = 284.625
$ws.Range("I7").Value = 224.5
$ws.Range("K7").Value = 224.5
$ws.Range("M7").Value = -111.5
$ws.Range("H11").Value = 403
$ws.Range("I11").Value = 403
$ws.Range("J11").Value = 0
$ws.Range("K11").Value = 403
$ws.Range("L11").Value = 0
$ws.Range("M11").Value = -263
$ws.Range("N11").ClearContents()
$ws.Range("H16").Value = 3000
$ws.Range("I16").Value = 3000
$ws.Range("K16").Value = 3000
$ws.Range("M16").Value = -2713
$ws.Range("H22").Value = 500
$ws.Range("J22").Value = 500
$ws.Range("L22").Value = 500
$ws.Range("N22").Value = -1200
$ws.Range("H95").Value = 40000
$ws.Range("J95").Value = 40000
$ws.Range("L95").Value = 40000
$ws.Range("N95").Value = -45492
$ws.Range("H107").Value = 677
$ws.Range("I107").Value = 689
$ws.Range("K107").Value = 689
$ws.Range("M107").Value = 1231
$ws.Range("H113").Value = 3000
$ws.Range("I113").Value = 3000
$ws.Range("K113").Value = 3000
$ws.Range("M113").Value = -830
$ws.Range("H135").Value = 15000
$ws.Range("J135").Value = 15000
$ws.Range("L135").Value = 15000
$ws.Range("N135").Value = -25140

$ws = $wb.Worksheets.Item("CUL")

$ws.Range("H2").Value = 41.583332
$ws.Range("I2").Value = 22.125
$ws.Range("J2").Value = 80.5
$ws.Range("K2").Value = 132.75
$ws.Range("L2").Value = 483
$ws.Range("M2").Value = -19.75
$ws.Range("N2").Value = -709
$ws.Range("H7").Value = 583.3333
$ws.Range("I7").Value = 100
$ws.Range("J7").Value = 3000
$ws.Range("K7").Value = 300
$ws.Range("L7").Value = 9000
$ws.Range("M7").Value = -188
$ws.Range("N7").Value = -9224
$ws.Range("H29").Value = 0
$ws.Range("I29").Value = 0
$ws.Range("K29").Value = 0
$ws.Range("M29").ClearContents()
$ws.Range("H50").Value = 0
$ws.Range("I50").Value = 0
$ws.Range("K50").Value = 0
$ws.Range("M50").ClearContents()
$ws.Range("H53").Value = 0
$ws.Range("I53").Value = 0
$ws.Range("K53").Value = 0
$ws.Range("M53").ClearContents()

$ws = $wb.Worksheets.Item("GSM")

$ws.Range("H21").Value = 666826.7
$ws.Range("I21").Value = 1000140
$ws.Range("J21").Value = 200
$ws.Range("K21").Value = 1000140
$ws.Range("L21").Value = 200
$ws.Range("M21").Value = -999967
$ws.Range("N21").Value = -546
$ws.Range("H30").Value = 666826.7
$ws.Range("I30").Value = 1000140
$ws.Range("J30").Value = 200
$ws.Range("K30").Value = 1000140
$ws.Range("L30").Value = 200
$ws.Range("M30").Value = -1000035
$ws.Range("N30").Value = -410
$ws.Range("H98").Value = 4800
$ws.Range("J98").Value = 4800
$ws.Range("L98").Value = 4800
$ws.Range("N98").Value = -10790
$ws.Range("H122").Value = 13080
$ws.Range("I122").Value = 16956.666
$ws.Range("J122").Value = 1450
$ws.Range("K122").Value = 50869.99800000001
$ws.Range("L122").Value = 4350
$ws.Range("M122").Value = -48419.99800000001
$ws.Range("N122").Value = -9250
$ws.Range("H126").Value = 4933.3335
$ws.Range("I126").Value = 4933.3335
$ws.Range("K126").Value = 14800.0005
$ws.Range("M126").Value = -12330.0005

$ws = $wb.Worksheets.Item("LTW")

$ws.Range("H16").Value = 1643.7778
$ws.Range("I16").Value = 1643.7778
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 1643.7778
$ws.Range("L16").Value = 0
$ws.Range("M16").Value = -1473.7778
$ws.Range("N16").ClearContents()
$ws.Range("H22").Value = 0
$ws.Range("J22").Value = 0
$ws.Range("L22").Value = 0
$ws.Range("N22").ClearContents()
$ws.Range("H27").Value = 0
$ws.Range("J27").Value = 0
$ws.Range("L27").Value = 0
$ws.Range("N27").ClearContents()
$ws.Range("H68").Value = 2135
$ws.Range("I68").Value = 1981.25
$ws.Range("K68").Value = 1981.25
$ws.Range("M68").Value = -1232.25
$ws.Range("H71").Value = 2135
$ws.Range("I71").Value = 1981.25
$ws.Range("K71").Value = 9906.25
$ws.Range("M71").Value = -6162.25
$ws.Range("H100").Value = 6597.1
$ws.Range("I100").Value = 2495.8572
$ws.Range("K100").Value = 2495.8572
$ws.Range("M100").Value = -1954.8572

$ws = $wb.Worksheets.Item("WVR")

$ws.Range("H2").Value = 874.38464
$ws.Range("I2").Value = 107.2
$ws.Range("J2").Value = 3431.6667
$ws.Range("K2").Value = 107.2
$ws.Range("L2").Value = 3431.6667
$ws.Range("M2").Value = 4.799999999999997
$ws.Range("N2").Value = -3655.6667
$ws.Range("H23").Value = 0
$ws.Range("I23").Value = 0
$ws.Range("J23").Value = 0
$ws.Range("K23").Value = 0
$ws.Range("L23").Value = 0
$ws.Range("M23").ClearContents()
$ws.Range("N23").ClearContents()

